# Update metadata obtained on 2016-04-06:
#  - Three columns (C, AY, BN) describing "sector-descripcion", "ratios" and
#    "rama-descripcion" are reclassified from iaest-measure to
#    iaest-dimension, their "medida" marker becomes "dim", and their data
#    type becomes "skos:Concept" instead of "xsd:string".
#  - A new metadata row (row 6) is added with the mapping file used for each
#    of those three dimensions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: iaest-measure:* -> iaest-dimension:* -------------------------
$ws.Range("C3").Value  = "iaest-dimension:sector-descripcion"
$ws.Range("AY3").Value = "iaest-dimension:ratios"
$ws.Range("BN3").Value = "iaest-dimension:rama-descripcion"

# --- Row 4: medida -> dim --------------------------------------------------
$ws.Range("C4").Value  = "dim"
$ws.Range("AY4").Value = "dim"
$ws.Range("BN4").Value = "dim"

# --- Row 5: xsd:string -> skos:Concept -------------------------------------
$ws.Range("C5").Value  = "skos:Concept"
$ws.Range("AY5").Value = "skos:Concept"
$ws.Range("BN5").Value = "skos:Concept"

# --- Row 6 (new): mapping file names ---------------------------------------
$ws.Range("C6").Value  = "mapping-sector-descripcion.xlsx"
$ws.Range("AY6").Value = "mapping-ratios.xlsx"
$ws.Range("BN6").Value = "mapping-rama-descripcion.xlsx"

# Match the formatting (style) used by the rest of the metadata rows.
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("AY5").Copy()
$ws.Range("AY6").PasteSpecial(-4122)
$ws.Range("BN5").Copy()
$ws.Range("BN6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
